$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.743.29"
$ws.Range("E2").Value = "  -1.18%  "
$ws.Range("D3").Value = "3.028.27"
$ws.Range("E3").Value = "  -1.33%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Formula = "'586.61"
$ws.Range("E5").Value = "  -0.70%  "
$ws.Range("D6").Formula = "'148.61"
$ws.Range("E6").Value = "  -3.59%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  -1.92%  "
$ws.Range("D9").Value = "3.025.95"
$ws.Range("E9").Value = "  -1.40%  "
$ws.Range("D10").Formula = "'0.151"
$ws.Range("E10").Value = "  -3.65%  "
$ws.Range("D11").Formula = "'5.88"
$ws.Range("E11").Value = "  -0.22%  "
$ws.Range("D12").Formula = "'0.452"
$ws.Range("E12").Value = "  +0.11%  "
$ws.Range("E13").Value = "  -2.84%  "
$ws.Range("D14").Formula = "'34.89"
$ws.Range("E14").Value = "  -5.03%  "
$ws.Range("E15").Value = "  +2.01%  "
$ws.Range("D16").Value = "3.524.37"
$ws.Range("E16").Value = "  -1.42%  "
$ws.Range("E17").Value = "  -0.61%  "
$ws.Range("D18").Value = "62.618.87"
$ws.Range("E18").Value = "  -1.26%  "
$ws.Range("D19").Value = "3.023.75"
$ws.Range("E19").Value = "  -1.36%  "
$ws.Range("D20").Formula = "'465.29"
$ws.Range("E20").Value = "  -3.78%  "
$ws.Range("D21").Formula = "'14.06"
$ws.Range("E21").Value = "  -3.62%  "
$ws.Range("D22").Formula = "'0.691"
$ws.Range("E22").Value = "  -2.59%  "
$ws.Range("D23").Formula = "'7.51"
$ws.Range("E23").Value = "  -0.76%  "
$ws.Range("E24").Value = "  -0.15%  "
$ws.Range("D25").Formula = "'2.28"
$ws.Range("E25").Value = "  -5.24%  "
$ws.Range("D26").Formula = "'12.45"
$ws.Range("E26").Value = "  -3.62%  "
$ws.Range("D27").Formula = "'10.42"
$ws.Range("E27").Value = "  -0.94%  "
$ws.Range("E28").Value = "  -0.12%  "
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("D30").Formula = "'2.64"
$ws.Range("E30").Value = "  -1.49%  "
$ws.Range("D31").Formula = "'7.19"
$ws.Range("E31").Value = "  -4.24%  "
$ws.Range("D32").Formula = "'2.13"
$ws.Range("E32").Value = "  -4.48%  "
$ws.Range("D33").Formula = "'28.75"
$ws.Range("E33").Value = "  +5.25%  "
$ws.Range("E34").Value = "  -2.33%  "
$ws.Range("D35").Value = "0.0₃0813"
$ws.Range("E35").Value = "  -1.45%  "
$ws.Range("E36").Value = "  -3.67%  "
$ws.Range("E37").Value = "  -4.06%  "
$ws.Range("D38").Formula = "'2.15"
$ws.Range("E38").Value = "  -3.52%  "
$ws.Range("D39").Formula = "'50.57"
$ws.Range("E39").Value = "  +0.08%  "
$ws.Range("D40").Formula = "'9.18"
$ws.Range("E40").Value = "  -1.60%  "
$ws.Range("D41").Formula = "'2.97"
$ws.Range("E41").Value = "  -8.66%  "
$ws.Range("E42").Value = "  +1.26%  "
$ws.Range("D43").Formula = "'400.10"
$ws.Range("E43").Value = "  -9.27%  "
$ws.Range("D44").Formula = "'0.279"
$ws.Range("E44").Value = "  -3.31%  "
$ws.Range("D45").Formula = "'0.0361"
$ws.Range("E45").Value = "  -0.88%  "
$ws.Range("D46").Value = "2.759.56"
$ws.Range("E46").Value = "  -2.41%  "
$ws.Range("D47").Formula = "'37.46"
$ws.Range("E47").Value = "  -5.33%  "
$ws.Range("D48").Formula = "'128.85"
$ws.Range("E48").Value = "  -3.00%  "
$ws.Range("E49").Value = "  +0.07%  "
$ws.Range("E50").Value = "  -0.25%  "
$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").Formula = "'2.21"
$ws.Range("E51").Value = "  -1.40%  "
